$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 39495
$ws.Range("E2").Value = 772518782873
$ws.Range("F2").Value = 15307010795
$ws.Range("G2").Value = 1.79578

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 2157.91
$ws.Range("E3").Value = 259558614486
$ws.Range("F3").Value = 16058033983
$ws.Range("G3").Value = 2.50008

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 1.001
$ws.Range("E4").Value = 89655985363
$ws.Range("F4").Value = 29490113955
$ws.Range("G4").Value = -0.00565

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 228.16
$ws.Range("E5").Value = 35109942764
$ws.Range("F5").Value = 566591965
$ws.Range("G5").Value = 0.15644

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "XRP"
$ws.Range("D6").Value = 0.626402
$ws.Range("E6").Value = 33742948711
$ws.Range("F6").Value = 775896413
$ws.Range("G6").Value = 1.98595

$ws.Range("B7").Value = "SOL"
$ws.Range("C7").Value = "Solana"
$ws.Range("D7").Value = 63.31
$ws.Range("E7").Value = 26884038459
$ws.Range("F7").Value = 1493532645
$ws.Range("G7").Value = 1.70548

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 24520889685
$ws.Range("F8").Value = 7462921120
$ws.Range("G8").Value = -0.01928

$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 2156.92
$ws.Range("E9").Value = 20025263480
$ws.Range("F9").Value = 12206808
$ws.Range("G9").Value = 2.70264

$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.39372
$ws.Range("E10").Value = 13779915710
$ws.Range("F10").Value = 281930639
$ws.Range("G10").Value = 1.05876

$ws.Range("B11").Value = "DOGE"
$ws.Range("C11").Value = "Dogecoin"
$ws.Range("D11").Value = 0.084936
$ws.Range("E11").Value = 12071076492
$ws.Range("F11").Value = 711440186
$ws.Range("G11").Value = 1.09001

$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.103278
$ws.Range("E12").Value = 9139206937
$ws.Range("F12").Value = 296902443
$ws.Range("G12").Value = 0.36781

$ws.Range("B13").Value = "LINK"
$ws.Range("C13").Value = "Chainlink"
$ws.Range("D13").Value = 16.04
$ws.Range("E13").Value = 8915592446
$ws.Range("F13").Value = 628175606
$ws.Range("G13").Value = 2.45334

$ws.Range("B14").Value = "AVAX"
$ws.Range("C14").Value = "Avalanche"
$ws.Range("D14").Value = 22
$ws.Range("E14").Value = 8037567030
$ws.Range("F14").Value = 330663452
$ws.Range("G14").Value = -0.26374

$ws.Range("B15").Value = "TON"
$ws.Range("C15").Value = "Toncoin"
$ws.Range("D15").Value = 2.32
$ws.Range("E15").Value = 8016239262
$ws.Range("F15").Value = 26926304
$ws.Range("G15").Value = -4.71039

$ws.Range("B16").Value = "MATIC"
$ws.Range("C16").Value = "Polygon"
$ws.Range("D16").Value = 0.81019
$ws.Range("E16").Value = 7515392716
$ws.Range("F16").Value = 469198570
$ws.Range("G16").Value = -0.29367

$ws.Range("B17").Value = "DOT"
$ws.Range("C17").Value = "Polkadot"
$ws.Range("D17").Value = 5.51
$ws.Range("E17").Value = 7178365396
$ws.Range("F17").Value = 123917949
$ws.Range("G17").Value = -0.36889

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 39439
$ws.Range("E18").Value = 6320098194
$ws.Range("F18").Value = 126027439
$ws.Range("G18").Value = 1.66095

$ws.Range("B19").Value = "LTC"
$ws.Range("C19").Value = "Litecoin"
$ws.Range("D19").Value = 72.3
$ws.Range("E19").Value = 5343801826
$ws.Range("F19").Value = 296141519
$ws.Range("G19").Value = 0.8087299999999999

$ws.Range("B20").Value = "DAI"
$ws.Range("C20").Value = "Dai"
$ws.Range("D20").Value = 0.998915
$ws.Range("E20").Value = 5327941206
$ws.Range("F20").Value = 314164095
$ws.Range("G20").Value = 0.03063

$ws.Range("B21").Value = "SHIB"
$ws.Range("C21").Value = "Shiba Inu"
$ws.Range("D21").Value = 0.00000848
$ws.Range("E21").Value = 4993193171
$ws.Range("F21").Value = 148225909
$ws.Range("G21").Value = 0.7681

$ws.Range("B22").Value = "UNI"
$ws.Range("C22").Value = "Uniswap"
$ws.Range("D22").Value = 6.18
$ws.Range("E22").Value = 4658480190
$ws.Range("F22").Value = 249085885
$ws.Range("G22").Value = 0.42131

$ws.Range("B23").Value = "BCH"
$ws.Range("C23").Value = "Bitcoin Cash"
$ws.Range("D23").Value = 227.9
$ws.Range("E23").Value = 4461502287
$ws.Range("F23").Value = 114129659
$ws.Range("G23").Value = 0.28814

$ws.Range("B24").Value = "OKB"
$ws.Range("C24").Value = "OKB"
$ws.Range("D24").Value = 58.08
$ws.Range("E24").Value = 3485569851
$ws.Range("F24").Value = 11062413
$ws.Range("G24").Value = -0.41796

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "LEO Token"
$ws.Range("D25").Value = 3.74
$ws.Range("E25").Value = 3471061565
$ws.Range("F25").Value = 1460084
$ws.Range("G25").Value = -6.02915

$ws.Range("B26").Value = "XLM"
$ws.Range("C26").Value = "Stellar"
$ws.Range("D26").Value = 0.121997
$ws.Range("E26").Value = 3419876393
$ws.Range("F26").Value = 51516454
$ws.Range("G26").Value = 1.66959

$ws.Range("B27").Value = "XMR"
$ws.Range("C27").Value = "Monero"
$ws.Range("D27").Value = 172.57
$ws.Range("E27").Value = 3131879300
$ws.Range("F27").Value = 61509738
$ws.Range("G27").Value = -0.06211

$ws.Range("B28").Value = "KAS"
$ws.Range("C28").Value = "Kaspa"
$ws.Range("D28").Value = 0.138882
$ws.Range("E28").Value = 3019095002
$ws.Range("F28").Value = 40140016
$ws.Range("G28").Value = 1.12971

$ws.Range("B29").Value = "TUSD"
$ws.Range("C29").Value = "TrueUSD"
$ws.Range("D29").Value = 0.9994690000000001
$ws.Range("E29").Value = 2954322170
$ws.Range("F29").Value = 143196884
$ws.Range("G29").Value = -0.04414

$ws.Range("B30").Value = "ETC"
$ws.Range("C30").Value = "Ethereum Classic"
$ws.Range("D30").Value = 19.62
$ws.Range("E30").Value = 2810403759
$ws.Range("F30").Value = 125222865
$ws.Range("G30").Value = 2.14564

$ws.Range("B31").Value = "ATOM"
$ws.Range("C31").Value = "Cosmos Hub"
$ws.Range("D31").Value = 9.52
$ws.Range("E31").Value = 2785997675
$ws.Range("F31").Value = 119433190
$ws.Range("G31").Value = -0.38629

$ws.Range("B32").Value = "CRO"
$ws.Range("C32").Value = "Cronos"
$ws.Range("D32").Value = 0.092071
$ws.Range("E32").Value = 2430527856
$ws.Range("F32").Value = 10056999
$ws.Range("G32").Value = 0.69543

$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "Filecoin"
$ws.Range("D33").Value = 4.61
$ws.Range("E33").Value = 2194137524
$ws.Range("F33").Value = 111528327
$ws.Range("G33").Value = 1.38422

$ws.Range("B34").Value = "LDO"
$ws.Range("C34").Value = "Lido DAO"
$ws.Range("D34").Value = 2.42
$ws.Range("E34").Value = 2161051507
$ws.Range("F34").Value = 61064034
$ws.Range("G34").Value = 1.07251

$ws.Range("B35").Value = "ICP"
$ws.Range("C35").Value = "Internet Computer"
$ws.Range("D35").Value = 4.7
$ws.Range("E35").Value = 2116239363
$ws.Range("F35").Value = 31640207
$ws.Range("G35").Value = -1.23783

$ws.Range("B36").Value = "RUNE"
$ws.Range("C36").Value = "THORChain"
$ws.Range("D36").Value = 6.98
$ws.Range("E36").Value = 2105945371
$ws.Range("F36").Value = 526307184
$ws.Range("G36").Value = -1.13234

$ws.Range("B37").Value = "HBAR"
$ws.Range("C37").Value = "Hedera"
$ws.Range("D37").Value = 0.062065
$ws.Range("E37").Value = 2083333046
$ws.Range("F37").Value = 23556789
$ws.Range("G37").Value = 0.54426

$ws.Range("B38").Value = "APT"
$ws.Range("C38").Value = "Aptos"
$ws.Range("D38").Value = 7.36
$ws.Range("E38").Value = 2046657856
$ws.Range("F38").Value = 87332896
$ws.Range("G38").Value = 1.91829

$ws.Range("B39").Value = "NEAR"
$ws.Range("C39").Value = "NEAR Protocol"
$ws.Range("D39").Value = 1.99
$ws.Range("E39").Value = 1992663764
$ws.Range("F39").Value = 136459021
$ws.Range("G39").Value = -1.2216

$ws.Range("B40").Value = "TAO"
$ws.Range("C40").Value = "Bittensor"
$ws.Range("D40").Value = 338.01
$ws.Range("E40").Value = 1946366825
$ws.Range("F40").Value = 13869063
$ws.Range("G40").Value = 17.26199

$ws.Range("B41").Value = "IMX"
$ws.Range("C41").Value = "Immutable"
$ws.Range("D41").Value = 1.42
$ws.Range("E41").Value = 1818049043
$ws.Range("F41").Value = 80758681
$ws.Range("G41").Value = 1.0308

$ws.Range("B42").Value = "MNT"
$ws.Range("C42").Value = "Mantle"
$ws.Range("D42").Value = 0.5425140000000001
$ws.Range("E42").Value = 1695251817
$ws.Range("F42").Value = 9834396
$ws.Range("G42").Value = -0.8942099999999999

$ws.Range("B43").Value = "BUSD"
$ws.Range("C43").Value = "BUSD"
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 1653773559
$ws.Range("F43").Value = 3266705450
$ws.Range("G43").Value = 0.00735

$ws.Range("B44").Value = "VET"
$ws.Range("C44").Value = "VeChain"
$ws.Range("D44").Value = 0.02255082
$ws.Range("E44").Value = 1640638062
$ws.Range("F44").Value = 42908257
$ws.Range("G44").Value = -0.28369

$ws.Range("B45").Value = "OP"
$ws.Range("C45").Value = "Optimism"
$ws.Range("D45").Value = 1.73
$ws.Range("E45").Value = 1570810069
$ws.Range("F45").Value = 113354116
$ws.Range("G45").Value = 1.37128

$ws.Range("B46").Value = "AAVE"
$ws.Range("C46").Value = "Aave"
$ws.Range("D46").Value = 102.35
$ws.Range("E46").Value = 1502030180
$ws.Range("F46").Value = 101653170
$ws.Range("G46").Value = -0.43851

$ws.Range("B47").Value = "INJ"
$ws.Range("C47").Value = "Injective"
$ws.Range("D47").Value = 17.81
$ws.Range("E47").Value = 1500172413
$ws.Range("F47").Value = 98274773
$ws.Range("G47").Value = -1.83675

$ws.Range("B48").Value = "QNT"
$ws.Range("C48").Value = "Quant"
$ws.Range("D48").Value = 100.7
$ws.Range("E48").Value = 1464908269
$ws.Range("F48").Value = 23279856
$ws.Range("G48").Value = 0.0106

$ws.Range("B49").Value = "ARB"
$ws.Range("C49").Value = "Arbitrum"
$ws.Range("D49").Value = 1.1
$ws.Range("E49").Value = 1404340942
$ws.Range("F49").Value = 264011806
$ws.Range("G49").Value = 2.29868

$ws.Range("B50").Value = "MKR"
$ws.Range("C50").Value = "Maker"
$ws.Range("D50").Value = 1524.39
$ws.Range("E50").Value = 1401150682
$ws.Range("F50").Value = 63798098
$ws.Range("G50").Value = -0.35539

$ws.Range("B51").Value = "GRT"
$ws.Range("C51").Value = "The Graph"
$ws.Range("D51").Value = 0.149269
$ws.Range("E51").Value = 1391613627
$ws.Range("F51").Value = 41776832
$ws.Range("G51").Value = -0.82706
